$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item(1)
$wsDetailed = $wb.Worksheets.Item(2)

# --- Schedule sheet (rows 3-5) ---
$wsSchedule.Range("E3").Value = 563.1449849999999
$wsSchedule.Range("F3").Value = 16.55335052910053
$wsSchedule.Range("B4").Value = 46038.10416666666
$wsSchedule.Range("C4").Value = 4.5
$wsSchedule.Range("D4").Value = 17.01
$wsSchedule.Range("E4").Value = 513.38060475
$wsSchedule.Range("F4").Value = 30.18110551146385
$wsSchedule.Range("A5").Value = 46038.27083333334
$wsSchedule.Range("C5").Value = 9.5
$wsSchedule.Range("D5").Value = 35.91
$wsSchedule.Range("E5").Value = 601.14417675
$wsSchedule.Range("F5").Value = 16.74030010442774

# --- Detailed sheet (rows 29-97) ---
$wsDetailed.Range("B29").Value = 36.0601
$wsDetailed.Range("B30").Value = 51.79253
$wsDetailed.Range("B31").Value = 39.7501
$wsDetailed.Range("C31").Value = "historical"
$wsDetailed.Range("B32").Value = 36.0601
$wsDetailed.Range("C32").Value = "historical"
$wsDetailed.Range("B33").Value = 57.04293
$wsDetailed.Range("C33").Value = "historical"
$wsDetailed.Range("B34").Value = 43.94456
$wsDetailed.Range("B35").Value = 46.29749
$wsDetailed.Range("B36").Value = 56.00261
$wsDetailed.Range("B37").Value = 24.10105
$wsDetailed.Range("B38").Value = 66.66182000000001
$wsDetailed.Range("B39").Value = 70.09057
$wsDetailed.Range("B40").Value = 100.58579
$wsDetailed.Range("B41").Value = 106.42876
$wsDetailed.Range("B42").Value = 102.84099
$wsDetailed.Range("B43").Value = 106.59
$wsDetailed.Range("B44").Value = 102.96758
$wsDetailed.Range("B45").Value = 85.95
$wsDetailed.Range("B46").Value = 67.89254
$wsDetailed.Range("B49").Value = 58.43705
$wsDetailed.Range("B50").Value = 57.09
$wsDetailed.Range("B51").Value = 57.92459
$wsDetailed.Range("B52").Value = 57.06003
$wsDetailed.Range("B53").Value = 56.98
$wsDetailed.Range("B54").Value = 56.98
$wsDetailed.Range("B55").Value = 57.06003
$wsDetailed.Range("E55").Value = "OFF"
$wsDetailed.Range("B56").Value = 56.98
$wsDetailed.Range("B57").Value = 56.88071
$wsDetailed.Range("B58").Value = 57.96129
$wsDetailed.Range("B59").Value = 59.63457
$wsDetailed.Range("B60").Value = 60.43567
$wsDetailed.Range("B61").Value = 65
$wsDetailed.Range("B62").Value = 58.27761
$wsDetailed.Range("E63").Value = "ON"
$wsDetailed.Range("B64").Value = 35.88
$wsDetailed.Range("B68").Value = 40.9658
$wsDetailed.Range("B69").Value = 36.06
$wsDetailed.Range("B71").Value = 40.7954
$wsDetailed.Range("B72").Value = 36.06028
$wsDetailed.Range("B73").Value = 36.06
$wsDetailed.Range("B77").Value = 27.4532
$wsDetailed.Range("B78").Value = 12.10384
$wsDetailed.Range("B79").Value = 7.36827
$wsDetailed.Range("B80").Value = 4.71079
$wsDetailed.Range("B81").Value = 29.70035
$wsDetailed.Range("B82").Value = 29.66805
$wsDetailed.Range("B83").Value = 24.42634
$wsDetailed.Range("B84").Value = 18.11114
$wsDetailed.Range("B85").Value = 16.09187
$wsDetailed.Range("B86").Value = 5.45427
$wsDetailed.Range("B87").Value = 40.36512
$wsDetailed.Range("B88").Value = 55.3303
$wsDetailed.Range("B89").Value = 67.39879999999999
$wsDetailed.Range("B90").Value = 59.37278
$wsDetailed.Range("B91").Value = 67.39879999999999
$wsDetailed.Range("B92").Value = 68.21745
$wsDetailed.Range("B93").Value = 57.04922
$wsDetailed.Range("B94").Value = 56.98078
$wsDetailed.Range("B95").Value = 48.38244
$wsDetailed.Range("B96").Value = 51.55592
$wsDetailed.Range("B97").Value = 57.02639
